$d = $word.ActiveDocument

# --- Step 1 -----------------------------------------------------------
# Collapse "My name is " + "himashu" + " " + "sarkar" + " and I know
# nothing" (with their spell-check proofErr wrappers) plus the trailing
# " of computer science and it." run into one clean sentence, in a
# single run. Using Find/Replace (rather than Range.Text=) also makes
# sure the edit is actually applied even though the visible characters
# end up identical to the original concatenated text.
$d.Content.Find.Execute(
    "My name is himashu sarkar and I know nothing of computer science and it.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "My name is himashu sarkar and I know nothing of computer science and it.",
    2) | Out-Null

# --- Step 2 -------------------------------------------------------------
# Add "what the fuck dude" as its own run. A plain InsertAfter right next
# to the previous run gets silently re-merged into it on save, so we
# temporarily put the new text in a throwaway paragraph of its own ...
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.End = $r.End - 1
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.End = $r2.End - 1
$r2.Collapse(0) | Out-Null
$r2.InsertAfter("what the fuck dude")

# ... then delete the paragraph mark that separated them again, so
# "what the fuck dude" rejoins paragraph 1 as a second run (it keeps its
# own run instead of merging with the first one because it was created
# in a different paragraph).
$p1 = $d.Paragraphs.Item(1)
$markRange = $d.Range($p1.Range.End - 1, $p1.Range.End)
$markRange.Delete()

# --- Step 3 ---------------------------------------------------------
# Insert a real paragraph break after "what the fuck dude" -- this new,
# empty paragraph is where the _GoBack bookmark ends up living on its
# own.
$p1 = $d.Paragraphs.Item(1)
$rEnd = $p1.Range
$rEnd.End = $rEnd.End - 1
$rEnd.Collapse(0) | Out-Null
$rEnd.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start

# --- Step 4 ---------------------------------------------------------
# Re-create the _GoBack bookmark inside the new, empty paragraph. Typing
# a placeholder character, wrapping the bookmark around it and then
# deleting the placeholder through the bookmark's own range leaves a
# clean zero-length bookmark (bookmarkStart immediately followed by
# bookmarkEnd, no stray empty run).
$tmp = $d.Range($p2Start, $p2Start)
$tmp.InsertAfter("X")
$bmRange = $d.Range($p2Start, $p2Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$bmr = $d.Bookmarks.Item("_GoBack").Range
$bmr.Text = ""
